$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1895.4166
$ws.Range("J17").Value = 1895.4166
$ws.Range("L17").Value = 5686.2498
$ws.Range("N17").Value = -6022.2498
$ws.Range("H33").Value = 333539.2
$ws.Range("J33").Value = 99.666664
$ws.Range("L33").Value = 99.666664
$ws.Range("N33").Value = -557.666664
$ws.Range("H121").Value = 3474.5
$ws.Range("J121").Value = 3474.5
$ws.Range("L121").Value = 10423.5
$ws.Range("N121").Value = -13917.5
$ws.Range("H125").Value = 3313.8572
$ws.Range("J125").Value = 4259.4
$ws.Range("L125").Value = 38334.6
$ws.Range("N125").Value = -43254.6
$ws.Range("H132").Value = 11180.482
$ws.Range("I132").Value = 3458.5
$ws.Range("K132").Value = 10375.5
$ws.Range("M132").Value = -7845.5
$ws.Range("H135").Value = 2137.5557
$ws.Range("I135").Value = 601.2857
$ws.Range("K135").Value = 5411.571300000001
$ws.Range("M135").Value = -2876.571300000001
$ws.Range("H138").Value = 2578
$ws.Range("I138").Value = 1629.4
$ws.Range("K138").Value = 4888.200000000001
$ws.Range("M138").Value = 251.7999999999993
$ws.Range("H141").Value = 7229.737
$ws.Range("I141").Value = 5773.6875
$ws.Range("J141").Value = 14995.333
$ws.Range("K141").Value = 17321.0625
$ws.Range("L141").Value = 44985.999
$ws.Range("M141").Value = -12141.0625
$ws.Range("N141").Value = -55345.999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 19235
$ws.Range("I61").Value = 16229.8
$ws.Range("K61").Value = 16229.8
$ws.Range("M61").Value = -16017.8
$ws.Range("H74").Value = 11364930
$ws.Range("I74").Value = 15626071
$ws.Range("J74").Value = 1887.8334
$ws.Range("K74").Value = 15626071
$ws.Range("L74").Value = 1887.8334
$ws.Range("M74").Value = -15625197
$ws.Range("N74").Value = -3635.8334
$ws.Range("H77").Value = 11364930
$ws.Range("I77").Value = 15626071
$ws.Range("J77").Value = 1887.8334
$ws.Range("K77").Value = 78130355
$ws.Range("L77").Value = 9439.166999999999
$ws.Range("M77").Value = -78125987
$ws.Range("N77").Value = -18175.167
$ws.Range("H122").Value = 2800
$ws.Range("I122").Value = 2800
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8400
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5950
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 6288.75
$ws.Range("I132").Value = 2275.2144
$ws.Range("J132").Value = 15653.667
$ws.Range("K132").Value = 6825.6432
$ws.Range("L132").Value = 46961.001
$ws.Range("M132").Value = -4295.6432
$ws.Range("N132").Value = -52021.001
$ws.Range("H136").Value = 19235
$ws.Range("I136").Value = 16229.8
$ws.Range("K136").Value = 48689.39999999999
$ws.Range("M136").Value = -46139.39999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3067.8635
$ws.Range("I105").Value = 2848.3
$ws.Range("J105").Value = 3250.8333
$ws.Range("K105").Value = 2848.3
$ws.Range("L105").Value = 3250.8333
$ws.Range("M105").Value = -1101.3
$ws.Range("N105").Value = -6744.8333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1446.8518
$ws.Range("I16").Value = 1577.6111
$ws.Range("J16").Value = 1185.3334
$ws.Range("K16").Value = 1577.6111
$ws.Range("L16").Value = 1185.3334
$ws.Range("M16").Value = -1290.6111
$ws.Range("N16").Value = -1759.3334
$ws.Range("H35").Value = 5525
$ws.Range("I35").Value = 5525
$ws.Range("K35").Value = 5525
$ws.Range("M35").Value = -5231
$ws.Range("H58").Value = 478645.72
$ws.Range("I58").Value = 716513.3
$ws.Range("K58").Value = 716513.3
$ws.Range("M58").Value = -716310.3
$ws.Range("H94").Value = 3038.875
$ws.Range("I94").Value = 2513
$ws.Range("K94").Value = 2513
$ws.Range("M94").Value = -2062
$ws.Range("H113").Value = 1446.8518
$ws.Range("I113").Value = 1577.6111
$ws.Range("J113").Value = 1185.3334
$ws.Range("K113").Value = 1577.6111
$ws.Range("L113").Value = 1185.3334
$ws.Range("M113").Value = 592.3888999999999
$ws.Range("N113").Value = -5525.3334
$ws.Range("H132").Value = 10754891
$ws.Range("I132").Value = 14494713
$ws.Range("K132").Value = 43484139
$ws.Range("M132").Value = -43481609
$ws.Range("H134").Value = 2410.2708
$ws.Range("I134").Value = 2527.641
$ws.Range("J134").Value = 1901.6666
$ws.Range("K134").Value = 7582.923000000001
$ws.Range("L134").Value = 5704.9998
$ws.Range("M134").Value = -5047.923000000001
$ws.Range("N134").Value = -10774.9998
$ws.Range("H136").Value = 478645.72
$ws.Range("I136").Value = 716513.3
$ws.Range("K136").Value = 2149539.9
$ws.Range("M136").Value = -2146989.9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 13570.286
$ws.Range("J112").Value = 14998.6
$ws.Range("L112").Value = 44995.8
$ws.Range("N112").Value = -47211.8
$ws.Range("H129").Value = 1520.2727
$ws.Range("I129").Value = 1823
$ws.Range("J129").Value = 1347.2858
$ws.Range("K129").Value = 5469
$ws.Range("L129").Value = 4041.8574
$ws.Range("M129").Value = -469
$ws.Range("N129").Value = -14041.8574
$ws.Range("H137").Value = 6484937
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 6484937
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 19454811
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -19465011
$ws.Range("H138").Value = 4015.3635
$ws.Range("I138").Value = 4024.2856
$ws.Range("J138").Value = 3999.75
$ws.Range("K138").Value = 12072.8568
$ws.Range("L138").Value = 11999.25
$ws.Range("M138").Value = -6932.856800000001
$ws.Range("N138").Value = -22279.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 19236732
$ws.Range("I102").Value = 26323140
$ws.Range("K102").Value = 26323140
$ws.Range("M102").Value = -26321518
$ws.Range("H122").Value = 788499.9
$ws.Range("I122").Value = 788499.9
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2365499.7
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2363049.7
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2887.1316
$ws.Range("I126").Value = 2697.7917
$ws.Range("J126").Value = 3211.7144
$ws.Range("K126").Value = 8093.375100000001
$ws.Range("L126").Value = 9635.143199999999
$ws.Range("M126").Value = -5623.375100000001
$ws.Range("N126").Value = -14575.1432
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2111.75
$ws.Range("I61").Value = 2081.5
$ws.Range("J61").Value = 2202.5
$ws.Range("K61").Value = 2081.5
$ws.Range("L61").Value = 2202.5
$ws.Range("M61").Value = -1879.5
$ws.Range("N61").Value = -2606.5
$ws.Range("H113").Value = 2111.75
$ws.Range("I113").Value = 2081.5
$ws.Range("J113").Value = 2202.5
$ws.Range("K113").Value = 2081.5
$ws.Range("L113").Value = 2202.5
$ws.Range("M113").Value = 88.5
$ws.Range("N113").Value = -6542.5
$ws.Range("H132").Value = 3401.84
$ws.Range("I132").Value = 2583.5667
$ws.Range("J132").Value = 6674.933
$ws.Range("K132").Value = 7750.7001
$ws.Range("L132").Value = 20024.799
$ws.Range("M132").Value = -5220.7001
$ws.Range("N132").Value = -25084.799
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 50000
$ws.Range("I75").Value = 50000
$ws.Range("K75").Value = 50000
$ws.Range("M75").Value = -49064
$ws.Range("H78").Value = 50000
$ws.Range("I78").Value = 50000
$ws.Range("K78").Value = 150000
$ws.Range("M78").Value = -145320
